# Generate Report for Handoff
# Adds a new row (for file "cac8ec24-cf23-4275-90be-d808920bb0a4ooo...md") to the
# Overview, zh-cn and de-de sheets/tables of the localization-status workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Shared literal values (kept identical across sheets)
# ---------------------------------------------------------------------------
$mdNew      = "cac8ec24-cf23-4275-90be-d808920bb0a4ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$e2eMdNew   = "e2e\cac8ec24-cf23-4275-90be-d808920bb0a4ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$readyForHandoff = "Ready for handoff"
$dt1        = "2016-08-25 10:29:06"
$zhXlfNew   = "cac8ec24-cf23-4275-90be-d808920bb0a4oooooooooooooooooooooooooooooooooooooooo.b2d4b67e2999042678ed599738f5d8b93074348b.zh-cn.xlf"
$dt2        = "2016-08-25 10:28:57"
$deXlfNew   = "cac8ec24-cf23-4275-90be-d808920bb0a4oooooooooooooooooooooooooooooooooooooooo.b2d4b67e2999042678ed599738f5d8b93074348b.de-de.xlf"
$hyperlinkUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/035d0049eb956264c4b61f02668c767471fdfd39/e2e/" + $mdNew

$linkUnderline = 2        # xlUnderlineStyleSingle
$linkColor     = 15570276 # BGR value of RGB(0x64,0x95,0xED) - matches existing HyperLink style
$dateFormat    = "yyyy-mm-dd HH:mm:ss"
$newColWidth   = 16.3      # closest achievable ColumnWidth to the target 17.216 characters

# ---------------------------------------------------------------------------
# Sheet "Overview" (table3 / sheet1) -- columns A:G
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A3").Value = $mdNew
$wsOverview.Range("B3").Value = $e2eMdNew
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("D3").Value = ""
$wsOverview.Range("E3").Value = $readyForHandoff
$wsOverview.Range("F3").Value = $readyForHandoff
$wsOverview.Range("G3").Value = $dt1
$wsOverview.Range("G3").NumberFormat = $dateFormat

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $hyperlinkUrl, "", "", $e2eMdNew) | Out-Null
$wsOverview.Range("B3").Font.Underline = $linkUnderline
$wsOverview.Range("B3").Font.Color = $linkColor

$wsOverview.Columns.Item(5).ColumnWidth = $newColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newColWidth

# ---------------------------------------------------------------------------
# Sheet "zh-cn" (table1 / sheet2) -- columns A:P
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$loZh.ListRows.Add() | Out-Null

$wsZh.Range("A3").Value = $mdNew
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = $readyForHandoff
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
$wsZh.Range("F3").Value = "False"
$wsZh.Range("G3").Value = $zhXlfNew
$wsZh.Range("H3").Value = $dt2
$wsZh.Range("H3").NumberFormat = $dateFormat
$wsZh.Range("I3").Value = ""
$wsZh.Range("J3").Value = ""
$wsZh.Range("K3").Value = "0001-01-01 00:00:00"
$wsZh.Range("K3").NumberFormat = $dateFormat
$wsZh.Range("L3").Value = ""
$wsZh.Range("M3").Value = "True"
$wsZh.Range("N3").Value = ""
$wsZh.Range("O3").Value = "False"
$wsZh.Range("P3").Value = ""

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $hyperlinkUrl, "", "", $mdNew) | Out-Null
$wsZh.Range("A3").Font.Underline = $linkUnderline
$wsZh.Range("A3").Font.Color = $linkColor

$wsZh.Columns.Item(3).ColumnWidth = $newColWidth

# ---------------------------------------------------------------------------
# Sheet "de-de" (table2 / sheet3) -- columns A:P
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$loDe.ListRows.Add() | Out-Null

$wsDe.Range("A3").Value = $mdNew
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = $readyForHandoff
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
$wsDe.Range("F3").Value = "False"
$wsDe.Range("G3").Value = $deXlfNew
$wsDe.Range("H3").Value = $dt1
$wsDe.Range("H3").NumberFormat = $dateFormat
$wsDe.Range("I3").Value = ""
$wsDe.Range("J3").Value = ""
$wsDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDe.Range("K3").NumberFormat = $dateFormat
$wsDe.Range("L3").Value = ""
$wsDe.Range("M3").Value = "True"
$wsDe.Range("N3").Value = ""
$wsDe.Range("O3").Value = "False"
$wsDe.Range("P3").Value = ""

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $hyperlinkUrl, "", "", $mdNew) | Out-Null
$wsDe.Range("A3").Font.Underline = $linkUnderline
$wsDe.Range("A3").Font.Color = $linkColor

$wsDe.Columns.Item(3).ColumnWidth = $newColWidth
